$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1248.7273
$ws.Range("I32").Value = 906.5
$ws.Range("K32").Value = 906.5
$ws.Range("M32").Value = -580.5

$ws.Range("H70").Value = 16730
$ws.Range("I70").Value = 1397
$ws.Range("J70").Value = 18433.666
$ws.Range("K70").Value = 4191
$ws.Range("L70").Value = 55300.99800000001
$ws.Range("M70").Value = -3921
$ws.Range("N70").Value = -55840.99800000001

$ws.Range("H73").Value = 16730
$ws.Range("I73").Value = 1397
$ws.Range("J73").Value = 18433.666
$ws.Range("K73").Value = 4191
$ws.Range("L73").Value = 55300.99800000001
$ws.Range("M73").Value = -3255
$ws.Range("N73").Value = -57172.99800000001

$ws.Range("H74").Value = 7651.2607
$ws.Range("I74").Value = 5230.1665
$ws.Range("K74").Value = 5230.1665
$ws.Range("M74").Value = -4294.1665

$ws.Range("H76").Value = 6859.533
$ws.Range("I76").Value = 5875.375
$ws.Range("J76").Value = 7984.2856
$ws.Range("K76").Value = 5875.375
$ws.Range("L76").Value = 7984.2856
$ws.Range("M76").Value = -5560.375
$ws.Range("N76").Value = -8614.285599999999

$ws.Range("H77").Value = 7651.2607
$ws.Range("I77").Value = 5230.1665
$ws.Range("K77").Value = 26150.8325
$ws.Range("M77").Value = -21470.8325

$ws.Range("H79").Value = 6859.533
$ws.Range("I79").Value = 5875.375
$ws.Range("J79").Value = 7984.2856
$ws.Range("K79").Value = 5875.375
$ws.Range("L79").Value = 7984.2856
$ws.Range("M79").Value = -4783.375
$ws.Range("N79").Value = -10168.2856

$ws.Range("H98").Value = 1275.8431
$ws.Range("I98").Value = 1085.7954
$ws.Range("K98").Value = 1085.7954
$ws.Range("M98").Value = 412.2046

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 0

$ws.Range("H122").Value = 1275.8431
$ws.Range("I122").Value = 1085.7954
$ws.Range("K122").Value = 3257.3862
$ws.Range("M122").Value = -807.3861999999999

$ws.Range("H131").Value = 5436.684
$ws.Range("I131").Value = 4878.357
$ws.Range("J131").Value = 7000
$ws.Range("K131").Value = 14635.071
$ws.Range("L131").Value = 21000
$ws.Range("M131").Value = -9595.071
$ws.Range("N131").Value = -31080

$ws.Range("H136").Value = 72500
$ws.Range("J136").Value = 72500
$ws.Range("L136").Value = 72500
$ws.Range("N136").Value = -82700

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 5132.1665
$ws.Range("J29").Value = 973.75
$ws.Range("L29").Value = 973.75
$ws.Range("N29").Value = -1589.75

$ws.Range("H32").Value = 5496733.5
$ws.Range("I32").Value = 5496733.5
$ws.Range("K32").Value = 5496733.5
$ws.Range("M32").Value = -5496446.5

$ws.Range("H45").Value = 3850.1428
$ws.Range("I45").Value = 3774.647
$ws.Range("K45").Value = 3774.647
$ws.Range("M45").Value = -3397.647

$ws.Range("H74").Value = 5438332.5
$ws.Range("I74").Value = 7355488.5
$ws.Range("J74").Value = 6390
$ws.Range("K74").Value = 7355488.5
$ws.Range("L74").Value = 6390
$ws.Range("M74").Value = -7354614.5
$ws.Range("N74").Value = -8138

$ws.Range("H77").Value = 5438332.5
$ws.Range("I77").Value = 7355488.5
$ws.Range("J77").Value = 6390
$ws.Range("K77").Value = 36777442.5
$ws.Range("L77").Value = 31950
$ws.Range("M77").Value = -36773074.5
$ws.Range("N77").Value = -40686

$ws.Range("H97").Value = 1190.7693
$ws.Range("I97").Value = 1289.6086
$ws.Range("K97").Value = 1289.6086
$ws.Range("M97").Value = -793.6086

$ws.Range("H132").Value = 973004.2
$ws.Range("I132").Value = 1165165.2
$ws.Range("K132").Value = 3495495.6
$ws.Range("M132").Value = -3492965.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2226.8667
$ws.Range("I86").Value = 2341.5557
$ws.Range("K86").Value = 2341.5557
$ws.Range("M86").Value = -1218.5557

$ws.Range("H89").Value = 2226.8667
$ws.Range("I89").Value = 2341.5557
$ws.Range("K89").Value = 11707.7785
$ws.Range("M89").Value = -6091.7785

$ws.Range("H94").Value = 1905.7142
$ws.Range("I94").Value = 1316.2142
$ws.Range("K94").Value = 1316.2142
$ws.Range("M94").Value = -865.2141999999999

$ws.Range("H134").Value = 805266.5600000001
$ws.Range("I134").Value = 1129835
$ws.Range("K134").Value = 3389505
$ws.Range("M134").Value = -3386970

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 240102.44
$ws.Range("I14").Value = 240102.44
$ws.Range("K14").Value = 720307.3200000001
$ws.Range("M14").Value = -720134.3200000001

$ws.Range("H107").Value = 886.7143
$ws.Range("I107").Value = 491.4
$ws.Range("J107").Value = 1875
$ws.Range("K107").Value = 1474.2
$ws.Range("L107").Value = 5625
$ws.Range("M107").Value = 445.8000000000002
$ws.Range("N107").Value = -9465

$ws.Range("H113").Value = 1994.375
$ws.Range("I113").Value = 1934.6666
$ws.Range("J113").Value = 2030.2
$ws.Range("K113").Value = 5803.9998
$ws.Range("L113").Value = 6090.6
$ws.Range("M113").Value = -3633.9998
$ws.Range("N113").Value = -10430.6

$ws.Range("H114").Value = 1721.8096
$ws.Range("I114").Value = 127.125
$ws.Range("J114").Value = 6824.8
$ws.Range("K114").Value = 381.375
$ws.Range("L114").Value = 20474.4
$ws.Range("M114").Value = 2872.625
$ws.Range("N114").Value = -26982.4

$ws.Range("H117").Value = 4822.2
$ws.Range("I117").Value = 2490.6
$ws.Range("J117").Value = 5988
$ws.Range("K117").Value = 7471.799999999999
$ws.Range("L117").Value = 17964
$ws.Range("M117").Value = -4029.799999999999
$ws.Range("N117").Value = -24848

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").ClearContents()
$ws.Range("N127").Value = 0

$ws.Range("H131").Value = 11099.739
$ws.Range("J131").Value = 16626.533
$ws.Range("L131").Value = 49879.599
$ws.Range("N131").Value = -59959.599

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3271.3333
$ws.Range("I122").Value = 3065.32
$ws.Range("K122").Value = 9195.960000000001
$ws.Range("M122").Value = -6745.960000000001

$ws.Range("H132").Value = 805294.0600000001
$ws.Range("I132").Value = 862386.5
$ws.Range("K132").Value = 2587159.5
$ws.Range("M132").Value = -2584629.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5692.25
$ws.Range("I40").Value = 4852.154
$ws.Range("J40").Value = 9332.666999999999
$ws.Range("K40").Value = 4852.154
$ws.Range("L40").Value = 9332.666999999999
$ws.Range("M40").Value = -4716.154
$ws.Range("N40").Value = -9604.666999999999

$ws.Range("H93").Value = 2281
$ws.Range("I93").Value = 2194.2
$ws.Range("J93").Value = 2329.2222
$ws.Range("K93").Value = 2194.2
$ws.Range("L93").Value = 2329.2222
$ws.Range("M93").Value = -946.1999999999998
$ws.Range("N93").Value = -4825.2222

$ws.Range("H122").Value = 3679.7317
$ws.Range("I122").Value = 3491.2693
$ws.Range("K122").Value = 10473.8079
$ws.Range("M122").Value = -8023.8079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1716.5
$ws.Range("I81").Value = 1654.8889
$ws.Range("J81").Value = 1901.3334
$ws.Range("K81").Value = 3309.7778
$ws.Range("L81").Value = 3802.6668
$ws.Range("M81").Value = -2248.7778
$ws.Range("N81").Value = -5924.6668

$ws.Range("H84").Value = 1716.5
$ws.Range("I84").Value = 1654.8889
$ws.Range("J84").Value = 1901.3334
$ws.Range("K84").Value = 16548.889
$ws.Range("L84").Value = 19013.334
$ws.Range("M84").Value = -11244.889
$ws.Range("N84").Value = -29621.334

$ws.Range("H96").Value = 2407.3845
$ws.Range("I96").Value = 2311.25
$ws.Range("J96").Value = 2561.2
$ws.Range("K96").Value = 2311.25
$ws.Range("L96").Value = 2561.2
$ws.Range("M96").Value = -938.25
$ws.Range("N96").Value = -5307.2

$ws.Range("H122").Value = 2142.7407
$ws.Range("I122").Value = 1617.75
$ws.Range("K122").Value = 4853.25
$ws.Range("M122").Value = -2403.25
